$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the wrongly-entered date in D92 (was off by exactly 100 years -> 76424, should be 43552)
$ws.Range("D92").Value = 43552

# Row 93: new journal entry - "Réalisation" activity about token management system
$ws.Range("A93").Value = "Réalisation"
$ws.Range("A93").WrapText = $true
$ws.Range("B93").Value = "Mise en place du nouveau système de gestion des token et incorporation de celui-ci sur l'interface web"
$ws.Range("B93").WrapText = $true
$ws.Range("C93").Value = 4.5
$ws.Range("D93").Value = 43553
$ws.Rows.Item(93).RowHeight = 60

# Row 94: new journal entry - "Documentation" activity about objectives summary
$ws.Range("A94").Value = "Documentation"
$ws.Range("A94").WrapText = $true
$ws.Range("B94").Value = "Définition des objectifs atteint et non-atteint. Détails des points m'ayant posé des problèmes particuliers"
$ws.Range("B94").WrapText = $true
$ws.Range("C94").Value = 1
$ws.Range("D94").Value = 43557
$ws.Rows.Item(94).RowHeight = 60

# Update the active cell selection to reflect where the user ended up working
$ws.Range("C95").Select()
